# fix typos in sample arc: C1_measured -> CC1_measured, C2_measured -> CC2_measured
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSEval")

$ws.Range("A2").Value = "CC1_measured"
$ws.Range("A3").Value = "CC2_measured"
$ws.Range("A4").Value = "CC2_measured"

# update playground.fsx / selection state: active cell moves from K7 to A4
$ws.Activate()
$ws.Range("A4").Select()
